$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45203 -> 45204) for every data row (rows 2 through 518).
$ws.Range("C2:C518").Value = 45204
